$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 422.9697
$ws.Range("I33").Value = 126.586205
$ws.Range("K33").Value = 126.586205
$ws.Range("M33").Value = 102.413795
$ws.Range("H86").Value = 9894
$ws.Range("I86").Value = 767.3333
$ws.Range("J86").Value = 20846
$ws.Range("K86").Value = 767.3333
$ws.Range("L86").Value = 20846
$ws.Range("M86").Value = 355.6667
$ws.Range("N86").Value = -23092
$ws.Range("H89").Value = 9894
$ws.Range("I89").Value = 767.3333
$ws.Range("J89").Value = 20846
$ws.Range("K89").Value = 3836.6665
$ws.Range("L89").Value = 104230
$ws.Range("M89").Value = 1779.3335
$ws.Range("N89").Value = -115462
$ws.Range("H98").Value = 827.9167
$ws.Range("I98").Value = 839.5454999999999
$ws.Range("K98").Value = 839.5454999999999
$ws.Range("M98").Value = 658.4545000000001
$ws.Range("H100").Value = 2680.7334
$ws.Range("I100").Value = 1617.5
$ws.Range("J100").Value = 3389.5557
$ws.Range("K100").Value = 1617.5
$ws.Range("L100").Value = 3389.5557
$ws.Range("M100").Value = -1076.5
$ws.Range("N100").Value = -4471.5557
$ws.Range("H106").Value = 9262161
$ws.Range("I106").Value = 23811270
$ws.Range("K106").Value = 23811270
$ws.Range("M106").Value = -23810639
$ws.Range("H122").Value = 827.9167
$ws.Range("I122").Value = 839.5454999999999
$ws.Range("K122").Value = 2518.6365
$ws.Range("M122").Value = -68.63649999999961
$ws.Range("H129").Value = 228163.89
$ws.Range("I129").Value = 263
$ws.Range("J129").Value = 286766.97
$ws.Range("K129").Value = 789
$ws.Range("L129").Value = 860300.9099999999
$ws.Range("M129").Value = 4211
$ws.Range("N129").Value = -870300.9099999999
$ws.Range("H137").Value = 120730.97
$ws.Range("I137").Value = 150431.6
$ws.Range("J137").Value = 6171.4287
$ws.Range("K137").Value = 451294.8
$ws.Range("L137").Value = 18514.2861
$ws.Range("M137").Value = -448744.8
$ws.Range("N137").Value = -23614.2861
$ws.Range("H138").Value = 4081.0305
$ws.Range("I138").Value = 3561.1875
$ws.Range("J138").Value = 4182.4634
$ws.Range("K138").Value = 10683.5625
$ws.Range("L138").Value = 12547.3902
$ws.Range("M138").Value = -5543.5625
$ws.Range("N138").Value = -22827.3902

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 26780
$ws.Range("J52").Value = 26780
$ws.Range("L52").Value = 26780
$ws.Range("N52").Value = -27416
$ws.Range("H61").Value = 11713992
$ws.Range("I61").Value = 15972480
$ws.Range("J61").Value = 3150
$ws.Range("K61").Value = 15972480
$ws.Range("L61").Value = 3150
$ws.Range("M61").Value = -15972268
$ws.Range("N61").Value = -3574
$ws.Range("H97").Value = 293.2
$ws.Range("I97").Value = 235.63637
$ws.Range("J97").Value = 451.5
$ws.Range("K97").Value = 235.63637
$ws.Range("L97").Value = 451.5
$ws.Range("M97").Value = 260.36363
$ws.Range("N97").Value = -1443.5
$ws.Range("H136").Value = 11713992
$ws.Range("I136").Value = 15972480
$ws.Range("J136").Value = 3150
$ws.Range("K136").Value = 47917440
$ws.Range("L136").Value = 9450
$ws.Range("M136").Value = -47914890
$ws.Range("N136").Value = -14550

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 22500
$ws.Range("I75").Value = 10000
$ws.Range("K75").Value = 10000
$ws.Range("M75").Value = -9064
$ws.Range("H78").Value = 22500
$ws.Range("I78").Value = 10000
$ws.Range("K78").Value = 30000
$ws.Range("M78").Value = -25320
$ws.Range("H134").Value = 5100
$ws.Range("I134").Value = 5296
$ws.Range("J134").Value = 4620.8887
$ws.Range("K134").Value = 15888
$ws.Range("L134").Value = 13862.6661
$ws.Range("M134").Value = -13353
$ws.Range("N134").Value = -18932.6661

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7139.564
$ws.Range("I31").Value = 3597.2
$ws.Range("K31").Value = 3597.2
$ws.Range("M31").Value = -3302.2
$ws.Range("H34").Value = 7139.564
$ws.Range("I34").Value = 3597.2
$ws.Range("K34").Value = 3597.2
$ws.Range("M34").Value = -3395.2
$ws.Range("H43").Value = 21999.5
$ws.Range("J43").Value = 21999.5
$ws.Range("L43").Value = 21999.5
$ws.Range("N43").Value = -22367.5
$ws.Range("H101").Value = 21999.5
$ws.Range("J101").Value = 21999.5
$ws.Range("L101").Value = 21999.5
$ws.Range("N101").Value = -28489.5
$ws.Range("H122").Value = 4837.6665
$ws.Range("I122").Value = 4837.6665
$ws.Range("K122").Value = 14512.9995
$ws.Range("M122").Value = -12062.9995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1324.9773
$ws.Range("I5").Value = 1132.3334
$ws.Range("J5").Value = 1630.9412
$ws.Range("K5").Value = 3397.0002
$ws.Range("L5").Value = 4892.8236
$ws.Range("M5").Value = -3285.0002
$ws.Range("N5").Value = -5116.8236
$ws.Range("H113").Value = 890
$ws.Range("J113").Value = 921.6667
$ws.Range("L113").Value = 2765.0001
$ws.Range("N113").Value = -7105.0001
$ws.Range("H122").Value = 1608.909
$ws.Range("J122").Value = 1661.7142
$ws.Range("L122").Value = 14955.4278
$ws.Range("N122").Value = -19855.4278
$ws.Range("H131").Value = 751.01
$ws.Range("J131").Value = 751.01
$ws.Range("L131").Value = 2253.03
$ws.Range("N131").Value = -12333.03
$ws.Range("H132").Value = 2257.1428
$ws.Range("I132").Value = 1033.3334
$ws.Range("J132").Value = 3175
$ws.Range("K132").Value = 9300.000599999999
$ws.Range("L132").Value = 28575
$ws.Range("M132").Value = -6770.000599999999
$ws.Range("N132").Value = -33635
$ws.Range("H133").Value = 6671.5386
$ws.Range("I133").Value = 2316.6667
$ws.Range("J133").Value = 7978
$ws.Range("K133").Value = 6950.000100000001
$ws.Range("L133").Value = 23934
$ws.Range("M133").Value = -1890.000100000001
$ws.Range("N133").Value = -34054
$ws.Range("H135").Value = 1324.9773
$ws.Range("I135").Value = 1132.3334
$ws.Range("J135").Value = 1630.9412
$ws.Range("K135").Value = 10191.0006
$ws.Range("L135").Value = 14678.4708
$ws.Range("M135").Value = -7656.000599999999
$ws.Range("N135").Value = -19748.4708

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3645.9583
$ws.Range("I80").Value = 3322.111
$ws.Range("J80").Value = 3840.2666
$ws.Range("K80").Value = 3322.111
$ws.Range("L80").Value = 3840.2666
$ws.Range("M80").Value = -2324.111
$ws.Range("N80").Value = -5836.2666
$ws.Range("H83").Value = 3645.9583
$ws.Range("I83").Value = 3322.111
$ws.Range("J83").Value = 3840.2666
$ws.Range("K83").Value = 16610.555
$ws.Range("L83").Value = 19201.333
$ws.Range("M83").Value = -11618.555
$ws.Range("N83").Value = -29185.333
$ws.Range("H126").Value = 5671.4287
$ws.Range("J126").Value = 5288.8887
$ws.Range("L126").Value = 15866.6661
$ws.Range("N126").Value = -20806.6661
$ws.Range("H132").Value = 9113887
$ws.Range("I132").Value = 25411802
$ws.Range("J132").Value = 59490.223
$ws.Range("K132").Value = 76235406
$ws.Range("L132").Value = 178470.669
$ws.Range("M132").Value = -76232876
$ws.Range("N132").Value = -183530.669
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3131.077
$ws.Range("I132").Value = 2188.889
$ws.Range("J132").Value = 5251
$ws.Range("K132").Value = 6566.667
$ws.Range("L132").Value = 15753
$ws.Range("M132").Value = -4036.667
$ws.Range("N132").Value = -20813
$ws.Range("H137").Value = 59273.8
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2783.077
$ws.Range("I126").Value = 2297.7778
$ws.Range("J126").Value = 3875
$ws.Range("K126").Value = 6893.3334
$ws.Range("L126").Value = 11625
$ws.Range("M126").Value = -4423.3334
$ws.Range("N126").Value = -16565
$ws.Range("H136").Value = 33337090
$ws.Range("I136").Value = 45456420
$ws.Range("K136").Value = 136369260
$ws.Range("M136").Value = -136366710
